$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("A12").Value = "05/08/2021 thru"
$ws.Range("B12").Value = "solving on pepecoding "
$ws.Range("C12").Value = "solved 1 question"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("G12").Value = "time and work"

# Row 13 - A13 reuses the date-number-format style from A9
$ws.Range("A9").Copy($ws.Range("A13"))
$ws.Range("A13").Value = "6/8/2021 Friday"
$ws.Range("B13").Value = "solving 1 problem "

# Update view: scroll/select like the author left it
$null = $ws.Range("J13").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
